# Weekly fruit/vegetable price update: a new weekly record is inserted as
# row 121 (pushing the existing rows 121-151 down to 122-152), matching the
# "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 121; Excel shifts rows 121:151 down to 122:152 and
# copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows.Item(121).Insert()

# Populate the newly inserted row 121 with this week's record.
$ws.Cells.Item(121, 1).Value  = 10
$ws.Cells.Item(121, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(121, 3).Value  = "La Araucanía"
$ws.Cells.Item(121, 4).Value  = 44722
$ws.Cells.Item(121, 5).Value  = 9
$ws.Cells.Item(121, 6).Value  = 100112012
$ws.Cells.Item(121, 7).Value  = "Espinaca"
$ws.Cells.Item(121, 8).Value  = "Sin especificar"
$ws.Cells.Item(121, 9).Value  = "Primera"
$ws.Cells.Item(121, 10).Value = 20
$ws.Cells.Item(121, 11).Value = 10000
$ws.Cells.Item(121, 12).Value = 10000
$ws.Cells.Item(121, 13).Value = 10000
$ws.Cells.Item(121, 14).Value = "$/docena de atados"
$ws.Cells.Item(121, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(121, 16).Value = 3333
$ws.Cells.Item(121, 17).Value = 3
$ws.Cells.Item(121, 18).Value = "Hortaliza"
